$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (rows 3-21): change the date separator from "/" to "-".
# A leading apostrophe forces Excel to keep these as text instead of
# auto-converting values such as "01-08-2022" into a date serial; the
# Style reset afterwards clears the quote-prefix formatting flag so no
# stray style gets attached to the cell.
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = "'" + $dates[$row]
    $cell.Style = "Normal"
}

# Update the Total Attendance Count (D), Real (E), Invalid (G) and Absent (H)
# counters for the rows whose values changed.
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 7).Value = 1

$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 8).Value = 0

$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 8).Value = 0

$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 8).Value = 0

$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 8).Value = 0

$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 8).Value = 0

$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 8).Value = 0

$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 8).Value = 0
